# QA Property Booking Bug list - "st4eps for debugging code explains flow"
# Updates Assignee (col C) values, clears some resolved Statuses (col H),
# removes the now-obsolete bug rows 15-20, and moves the sheet selection back
# to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear "Resolved" status that no longer applies ---------------------
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("H11").ClearContents()

# --- Reassign bugs to their new owners -----------------------------------
$ws.Range("C10").Value = "Shruti"
$ws.Range("C11").Value = "Shruti"
$ws.Range("C12").Value = "Shruti"
$ws.Range("C13").Value = "Mukesh/Shruti"
$ws.Range("C14").Value = "Mukesh/Shruti"

# --- Remove the stale/duplicate bug rows 15-20 ---------------------------
$ws.Rows("15:20").Delete() | Out-Null

# --- Reset the view: scroll to top, select I5 ----------------------------
$ws.Range("I5").Select() | Out-Null
